# Applies the commit "Updated symbol list on Wed Jan 25 20:42:54 UTC 2023
# with GitHub Actions" to the crypto price table on the active worksheet.
#
# The update refreshes the Price (column D) and Volume(1h) (column E)
# readings for most of the existing rows, and rotates a new coin
# (MXToken) into the top-20 block: every row from 8 to 17 shifts down by
# one position (GateToken -> MXToken -> LiechtensteinCryptoassetsExchange
# -> WazirX -> MandalaExchangeToken -> BitrueCoin -> BitMartToken ->
# BitForexToken -> TigerCash -> LEO -> GateToken), each carrying its own
# refreshed Coin name / Link / Price / Volume values along with it.
#
# All of the Coin/Link/Price/Volume cells on this sheet are stored as
# plain text (the sheet never uses real numeric or percentage cells), so
# every numeric-looking replacement value below is written with a leading
# apostrophe. That forces Excel to keep storing it as literal text
# instead of silently reinterpreting it as a float/percentage.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = "D2"; Value = "'303.28" },
    @{ Cell = "E2"; Value = "'-2.07%" },
    @{ Cell = "E3"; Value = "'0.35%" },
    @{ Cell = "E4"; Value = "'-0.74%" },
    @{ Cell = "D5"; Value = "'0.08069" },
    @{ Cell = "E5"; Value = "'-1.55%" },
    @{ Cell = "D6"; Value = "'1.920" },
    @{ Cell = "E6"; Value = "'-6.20%" },
    @{ Cell = "D7"; Value = "'7.789" },
    @{ Cell = "E7"; Value = "'-2.29%" },
    @{ Cell = "B8"; Value = "MXToken" },
    @{ Cell = "C8"; Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx" },
    @{ Cell = "D8"; Value = "'0.9292" },
    @{ Cell = "E8"; Value = "'0.04%" },
    @{ Cell = "B9"; Value = "LiechtensteinCryptoassetsExchange" },
    @{ Cell = "C9"; Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx" },
    @{ Cell = "D9"; Value = "'0.1493" },
    @{ Cell = "E9"; Value = "'37.68%" },
    @{ Cell = "B10"; Value = "WazirX" },
    @{ Cell = "C10"; Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx" },
    @{ Cell = "D10"; Value = "'0.1896" },
    @{ Cell = "E10"; Value = "'-1.58%" },
    @{ Cell = "B11"; Value = "MandalaExchangeToken" },
    @{ Cell = "C11"; Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx" },
    @{ Cell = "D11"; Value = "'0.09009" },
    @{ Cell = "E11"; Value = "'-5.51%" },
    @{ Cell = "B12"; Value = "BitrueCoin" },
    @{ Cell = "C12"; Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr" },
    @{ Cell = "D12"; Value = "'0.03444" },
    @{ Cell = "E12"; Value = "'-3.88%" },
    @{ Cell = "B13"; Value = "BitMartToken" },
    @{ Cell = "C13"; Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx" },
    @{ Cell = "D13"; Value = "'0.09843" },
    @{ Cell = "E13"; Value = "'-0.61%" },
    @{ Cell = "B14"; Value = "BitForexToken" },
    @{ Cell = "C14"; Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf" },
    @{ Cell = "D14"; Value = "'0.001418" },
    @{ Cell = "E14"; Value = "'-0.67%" },
    @{ Cell = "B15"; Value = "TigerCash" },
    @{ Cell = "C15"; Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch" },
    @{ Cell = "D15"; Value = "'0.005755" },
    @{ Cell = "E15"; Value = "'0.49%" },
    @{ Cell = "B16"; Value = "LEO" },
    @{ Cell = "C16"; Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo" },
    @{ Cell = "D16"; Value = "'3.545" },
    @{ Cell = "E16"; Value = "'2.19%" },
    @{ Cell = "B17"; Value = "GateToken" },
    @{ Cell = "C17"; Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt" },
    @{ Cell = "D17"; Value = "'4.052" },
    @{ Cell = "E17"; Value = "'-1.85%" },
    @{ Cell = "E18"; Value = "'3.00%" },
    @{ Cell = "D19"; Value = "'0.3441" },
    @{ Cell = "E19"; Value = "'-0.43%" },
    @{ Cell = "D20"; Value = "'0.1299" },
    @{ Cell = "E20"; Value = "'-0.89%" },
    @{ Cell = "D21"; Value = "'5.029" },
    @{ Cell = "E21"; Value = "'-1.43%" },
    @{ Cell = "D22"; Value = "'0.2492" },
    @{ Cell = "E22"; Value = "'13.69%" },
    @{ Cell = "D23"; Value = "'0.04504" },
    @{ Cell = "E23"; Value = "'-1.19%" },
    @{ Cell = "E24"; Value = "'-1.49%" },
    @{ Cell = "D25"; Value = "'0.004809" },
    @{ Cell = "E25"; Value = "'0.31%" },
    @{ Cell = "D26"; Value = "'0.0001227" },
    @{ Cell = "E26"; Value = "'-2.03%" },
    @{ Cell = "E27"; Value = "'-32.27%" },
    @{ Cell = "D39"; Value = "'0.01873" },
    @{ Cell = "E39"; Value = "'-5.34%" },
    @{ Cell = "D40"; Value = "'0.04788" },
    @{ Cell = "E40"; Value = "'-2.11%" },
    @{ Cell = "E41"; Value = "'7.52%" },
    @{ Cell = "D42"; Value = "'0.007326" },
    @{ Cell = "E42"; Value = "'-4.46%" },
    @{ Cell = "E43"; Value = "'-2.70%" },
    @{ Cell = "D44"; Value = "'0.002105" },
    @{ Cell = "E44"; Value = "'-0.67%" },
    @{ Cell = "D45"; Value = "'0.009728" },
    @{ Cell = "E45"; Value = "'-15.90%" },
    @{ Cell = "D46"; Value = "'0.00006212" },
    @{ Cell = "E46"; Value = "'-4.79%" },
    @{ Cell = "E47"; Value = "'-0.46%" },
    @{ Cell = "E48"; Value = "'-63.13%" },
    @{ Cell = "E50"; Value = "'-0.46%" },
    @{ Cell = "D51"; Value = "'0.0001995" },
    @{ Cell = "E51"; Value = "'-0.46%" }

)

foreach ($update in $updates) {
    $ws.Range($update.Cell).Value = $update.Value
}

